$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $text) {
    $origStyle = $cell.Style
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.Style = $origStyle
}

Set-TextValue $ws.Range('D2') '40.807.19'
Set-TextValue $ws.Range('E2') '  -2.17%  '

Set-TextValue $ws.Range('D3') '2.380.83'
Set-TextValue $ws.Range('E3') '  -3.73%  '

Set-TextValue $ws.Range('D4') '1.00'
Set-TextValue $ws.Range('E4') '  +0.15%  '

Set-TextValue $ws.Range('D5') '312.84'
Set-TextValue $ws.Range('E5') '  -1.97%  '

Set-TextValue $ws.Range('D6') '88.26'
Set-TextValue $ws.Range('E6') '  -5.40%  '

Set-TextValue $ws.Range('E7') '  -3.91%  '

Set-TextValue $ws.Range('E8') '  +0.02%  '

Set-TextValue $ws.Range('D9') '0.493'
Set-TextValue $ws.Range('E9') '  -4.72%  '

Set-TextValue $ws.Range('D10') '0.0835'
Set-TextValue $ws.Range('E10') '  -3.41%  '

Set-TextValue $ws.Range('D11') '31.15'
Set-TextValue $ws.Range('E11') '  -6.68%  '

Set-TextValue $ws.Range('E12') '  -1.46%  '

Set-TextValue $ws.Range('D13') '2.753.97'
Set-TextValue $ws.Range('E13') '  -3.54%  '

Set-TextValue $ws.Range('E14') '  -4.86%  '

Set-TextValue $ws.Range('D15') '15.09'
Set-TextValue $ws.Range('E15') '  -4.41%  '

Set-TextValue $ws.Range('D16') '2.377.34'
Set-TextValue $ws.Range('E16') '  -4.31%  '

Set-TextValue $ws.Range('D17') '0.764'
Set-TextValue $ws.Range('E17') '  -3.93%  '

Set-TextValue $ws.Range('D18') '40.760.34'
Set-TextValue $ws.Range('E18') '  -2.15%  '

Set-TextValue $ws.Range('D19') '0.0₃0916'
Set-TextValue $ws.Range('E19') '  -3.59%  '

Set-TextValue $ws.Range('D20') '6.19'
Set-TextValue $ws.Range('E20') '  -4.23%  '

Set-TextValue $ws.Range('D21') '69.28'
Set-TextValue $ws.Range('E21') '  -2.69%  '

Set-TextValue $ws.Range('D22') '10.79'
Set-TextValue $ws.Range('E22') '  -4.52%  '

Set-TextValue $ws.Range('D23') '233.61'
Set-TextValue $ws.Range('E23') '  -2.45%  '

Set-TextValue $ws.Range('E24') '  -3.37%  '

Set-TextValue $ws.Range('E25') '  +0.04%  '

Set-TextValue $ws.Range('D26') '1.82'
Set-TextValue $ws.Range('E26') '  -6.13%  '

Set-TextValue $ws.Range('D27') '23.85'
Set-TextValue $ws.Range('E27') '  -3.46%  '

Set-TextValue $ws.Range('D28') '2.20'
Set-TextValue $ws.Range('E28') '  -2.45%  '

Set-TextValue $ws.Range('D29') '9.39'
Set-TextValue $ws.Range('E29') '  -4.17%  '

Set-TextValue $ws.Range('D30') '33.86'
Set-TextValue $ws.Range('E30') '  -6.09%  '

Set-TextValue $ws.Range('D31') '154.79'
Set-TextValue $ws.Range('E31') '  -2.74%  '

Set-TextValue $ws.Range('E32') '  +0.14%  '

Set-TextValue $ws.Range('D33') '5.20'
Set-TextValue $ws.Range('E33') '  -5.66%  '

Set-TextValue $ws.Range('D34') '0.0737'
Set-TextValue $ws.Range('E34') '  -3.86%  '

Set-TextValue $ws.Range('E35') '  -6.13%  '

Set-TextValue $ws.Range('E36') '  -2.15%  '

Set-TextValue $ws.Range('D37') '2.81'
Set-TextValue $ws.Range('E37') '  -3.92%  '

Set-TextValue $ws.Range('D38') '16.06'
Set-TextValue $ws.Range('E38') '  -8.41%  '

Set-TextValue $ws.Range('D39') '0.0998'
Set-TextValue $ws.Range('E39') '  -3.52%  '

Set-TextValue $ws.Range('E40') '  -7.60%  '

Set-TextValue $ws.Range('D41') '3.81'
Set-TextValue $ws.Range('E41') '  -5.48%  '

Set-TextValue $ws.Range('D42') '2.33'
Set-TextValue $ws.Range('E42') '  -5.74%  '

Set-TextValue $ws.Range('D43') '1.961.92'
Set-TextValue $ws.Range('E43') '  -1.74%  '

Set-TextValue $ws.Range('D44') '0.0271'
Set-TextValue $ws.Range('E44') '  -5.03%  '

Set-TextValue $ws.Range('D45') '17.66'
Set-TextValue $ws.Range('E45') '  -6.52%  '

Set-TextValue $ws.Range('D46') '2.79'
Set-TextValue $ws.Range('E46') '  -7.02%  '

Set-TextValue $ws.Range('D47') '9.32'
Set-TextValue $ws.Range('E47') '  -1.61%  '

Set-TextValue $ws.Range('D48') '2.617.87'
Set-TextValue $ws.Range('E48') '  -3.46%  '

Set-TextValue $ws.Range('B49') 'BitcoinSV'
Set-TextValue $ws.Range('C49') 'https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv'
Set-TextValue $ws.Range('D49') '73.11'
Set-TextValue $ws.Range('E49') '  -1.12%  '

Set-TextValue $ws.Range('B50') 'Aave'
Set-TextValue $ws.Range('C50') 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
Set-TextValue $ws.Range('D50') '93.71'
Set-TextValue $ws.Range('E50') '  -3.79%  '

Set-TextValue $ws.Range('D51') '50.96'
Set-TextValue $ws.Range('E51') '  -3.37%  '
